$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 744; this shifts the existing rows 744-797
# down to 745-798 (matches dimension change A1:T797 -> A1:T798).
$ws.Range("A744:T744").EntireRow.Insert()

# Populate the newly inserted row 744 with the new record.
$ws.Cells.Item(744, 1).Value = 6
$ws.Cells.Item(744, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(744, 3).Value = "Metropolitana"
$ws.Cells.Item(744, 4).Value = 45021
$ws.Cells.Item(744, 5).Value = 13
$ws.Cells.Item(744, 6).Value = "Fruta"
$ws.Cells.Item(744, 7).Value = 100103
$ws.Cells.Item(744, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(744, 9).Value = 100103002
$ws.Cells.Item(744, 10).Value = "Ciruela"
$ws.Cells.Item(744, 11).Value = "Angeleno"
$ws.Cells.Item(744, 12).Value = "Segunda"
$ws.Cells.Item(744, 13).Value = 24
$ws.Cells.Item(744, 14).Value = 120000
$ws.Cells.Item(744, 15).Value = 130000
$ws.Cells.Item(744, 16).Value = 125000
$ws.Cells.Item(744, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(744, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(744, 19).Value = 278
$ws.Cells.Item(744, 20).Value = 450
